$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1179.4667
$ws.Range("I80").Value = 347
$ws.Range("J80").Value = 1482.1818
$ws.Range("K80").Value = 1041
$ws.Range("L80").Value = 4446.5454
$ws.Range("M80").Value = -43
$ws.Range("N80").Value = -6442.5454
$ws.Range("H83").Value = 1179.4667
$ws.Range("I83").Value = 347
$ws.Range("J83").Value = 1482.1818
$ws.Range("K83").Value = 3123
$ws.Range("L83").Value = 13339.6362
$ws.Range("M83").Value = 1869
$ws.Range("N83").Value = -23323.6362
$ws.Range("H86").Value = 4193.6
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 3992
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 3992
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -6238
$ws.Range("H88").Value = 3658.1875
$ws.Range("I88").Value = 766.8
$ws.Range("J88").Value = 4972.4546
$ws.Range("K88").Value = 766.8
$ws.Range("L88").Value = 4972.4546
$ws.Range("M88").Value = -360.8
$ws.Range("N88").Value = -5784.4546
$ws.Range("H89").Value = 4193.6
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 3992
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 19960
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -31192
$ws.Range("H91").Value = 3658.1875
$ws.Range("I91").Value = 766.8
$ws.Range("J91").Value = 4972.4546
$ws.Range("K91").Value = 766.8
$ws.Range("L91").Value = 4972.4546
$ws.Range("M91").Value = 637.2
$ws.Range("N91").Value = -7780.4546
$ws.Range("H129").Value = 1704.579
$ws.Range("J129").Value = 2078.1538
$ws.Range("L129").Value = 6234.4614
$ws.Range("N129").Value = -16234.4614

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5128.5713
$ws.Range("I74").Value = 5314.3335
$ws.Range("J74").Value = 4014
$ws.Range("K74").Value = 5314.3335
$ws.Range("L74").Value = 4014
$ws.Range("M74").Value = -4440.3335
$ws.Range("N74").Value = -5762
$ws.Range("H77").Value = 5128.5713
$ws.Range("I77").Value = 5314.3335
$ws.Range("J77").Value = 4014
$ws.Range("K77").Value = 26571.6675
$ws.Range("L77").Value = 20070
$ws.Range("M77").Value = -22203.6675
$ws.Range("N77").Value = -28806
$ws.Range("H88").Value = 3630
$ws.Range("I88").Value = 890
$ws.Range("K88").Value = 890
$ws.Range("M88").Value = -484
$ws.Range("H91").Value = 3630
$ws.Range("I91").Value = 890
$ws.Range("K91").Value = 890
$ws.Range("M91").Value = 514

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5625.5
$ws.Range("I94").Value = 4688.625
$ws.Range("J94").Value = 7499.25
$ws.Range("K94").Value = 4688.625
$ws.Range("L94").Value = 7499.25
$ws.Range("M94").Value = -4237.625
$ws.Range("N94").Value = -8401.25
$ws.Range("H99").Value = 3714.238
$ws.Range("I99").Value = 2566.6667
$ws.Range("K99").Value = 2566.6667
$ws.Range("M99").Value = -1068.6667
$ws.Range("H130").Value = 289999
$ws.Range("J130").Value = 289999
$ws.Range("L130").Value = 289999
$ws.Range("N130").Value = -300039
$ws.Range("H131").Value = 232110.11
$ws.Range("J131").Value = 227374
$ws.Range("L131").Value = 227374
$ws.Range("N131").Value = -237454

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2519.75
$ws.Range("J16").Value = 4999
$ws.Range("L16").Value = 4999
$ws.Range("N16").Value = -5573
$ws.Range("H31").Value = 3888.75
$ws.Range("I31").Value = 1731.2727
$ws.Range("J31").Value = 11799.5
$ws.Range("K31").Value = 1731.2727
$ws.Range("L31").Value = 11799.5
$ws.Range("M31").Value = -1436.2727
$ws.Range("N31").Value = -12389.5
$ws.Range("H34").Value = 3888.75
$ws.Range("I34").Value = 1731.2727
$ws.Range("J34").Value = 11799.5
$ws.Range("K34").Value = 1731.2727
$ws.Range("L34").Value = 11799.5
$ws.Range("M34").Value = -1529.2727
$ws.Range("N34").Value = -12203.5
$ws.Range("H58").Value = 2423.75
$ws.Range("I58").Value = 2423.75
$ws.Range("K58").Value = 2423.75
$ws.Range("M58").Value = -2220.75
$ws.Range("H113").Value = 2519.75
$ws.Range("J113").Value = 4999
$ws.Range("L113").Value = 4999
$ws.Range("N113").Value = -9339
$ws.Range("H132").Value = 2537.7273
$ws.Range("I132").Value = 2394.6072
$ws.Range("K132").Value = 7183.821599999999
$ws.Range("M132").Value = -4653.821599999999
$ws.Range("H136").Value = 2423.75
$ws.Range("I136").Value = 2423.75
$ws.Range("K136").Value = 7271.25
$ws.Range("M136").Value = -4721.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 15
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 15
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 45
$ws.Range("N75").Value = -2041
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 15
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 15
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 135
$ws.Range("N78").Value = -10119
$ws.Range("M78").ClearContents()
$ws.Range("H105").Value = 53343
$ws.Range("J105").Value = 53343
$ws.Range("L105").Value = 160029
$ws.Range("N105").Value = -165271

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 22975
$ws.Range("J15").Value = 23969
$ws.Range("L15").Value = 23969
$ws.Range("N15").Value = -24545
$ws.Range("H54").Value = 7000
$ws.Range("J54").Value = 7000
$ws.Range("L54").Value = 7000
$ws.Range("N54").Value = -7780
$ws.Range("H80").Value = 3655.3333
$ws.Range("I80").Value = 3632.6667
$ws.Range("J80").Value = 3666.6667
$ws.Range("K80").Value = 3632.6667
$ws.Range("L80").Value = 3666.6667
$ws.Range("M80").Value = -2634.6667
$ws.Range("N80").Value = -5662.6667
$ws.Range("H81").Value = 22975
$ws.Range("J81").Value = 23969
$ws.Range("L81").Value = 23969
$ws.Range("N81").Value = -25965
$ws.Range("H83").Value = 3655.3333
$ws.Range("I83").Value = 3632.6667
$ws.Range("J83").Value = 3666.6667
$ws.Range("K83").Value = 18163.3335
$ws.Range("L83").Value = 18333.3335
$ws.Range("M83").Value = -13171.3335
$ws.Range("N83").Value = -28317.3335
$ws.Range("H84").Value = 22975
$ws.Range("J84").Value = 23969
$ws.Range("L84").Value = 71907
$ws.Range("N84").Value = -81891
$ws.Range("H132").Value = 2102.1853
$ws.Range("I132").Value = 2181.2917
$ws.Range("J132").Value = 1469.3334
$ws.Range("K132").Value = 6543.875100000001
$ws.Range("L132").Value = 4408.0002
$ws.Range("M132").Value = -4013.875100000001
$ws.Range("N132").Value = -9468.0002

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3457.3157
$ws.Range("I22").Value = 3243.5
$ws.Range("J22").Value = 3612.818
$ws.Range("K22").Value = 3243.5
$ws.Range("L22").Value = 3612.818
$ws.Range("M22").Value = -2948.5
$ws.Range("N22").Value = -4202.818
$ws.Range("H27").Value = 3457.3157
$ws.Range("I27").Value = 3243.5
$ws.Range("J27").Value = 3612.818
$ws.Range("K27").Value = 3243.5
$ws.Range("L27").Value = 3612.818
$ws.Range("M27").Value = -3136.5
$ws.Range("N27").Value = -3826.818
$ws.Range("H93").Value = 4994.6665
$ws.Range("I93").Value = 2433.3333
$ws.Range("K93").Value = 2433.3333
$ws.Range("M93").Value = -1185.3333
$ws.Range("H100").Value = 6354
$ws.Range("I100").Value = 1166.3334
$ws.Range("J100").Value = 10244.75
$ws.Range("K100").Value = 1166.3334
$ws.Range("L100").Value = 10244.75
$ws.Range("M100").Value = -625.3334
$ws.Range("N100").Value = -11326.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1483
$ws.Range("I81").Value = 1541.9166
$ws.Range("J81").Value = 1306.25
$ws.Range("K81").Value = 3083.8332
$ws.Range("L81").Value = 2612.5
$ws.Range("M81").Value = -2022.8332
$ws.Range("N81").Value = -4734.5
$ws.Range("H84").Value = 1483
$ws.Range("I84").Value = 1541.9166
$ws.Range("J84").Value = 1306.25
$ws.Range("K84").Value = 15419.166
$ws.Range("L84").Value = 13062.5
$ws.Range("M84").Value = -10115.166
$ws.Range("N84").Value = -23670.5
